# "Maps 2 RGB" sheet: fill in the previously-blank "V4.07.5" .. "V4.07.8"
# result rows (29-32) with the measured PSNR/SSIM numbers, matching the
# formatting already used by the rows above them (B25:I28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps 2 RGB")

# Copy the number format (0.0000, default black font) from the row right
# above down onto the new rows -- these cells previously had an empty
# "pending" red-font style; this replaces it with the normal look used by
# every other completed row.
$ws.Range("B28:I28").Copy() | Out-Null
$ws.Range("B29:I32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$results = @(
    @(7.6207, 0.6047, 16.1814, 0.7642, 16.8018, 0.8069, 22.2549, 0.8223),
    @(7.6005, 0.5768, 12.7917, 0.6454, 15.7297, 0.7965, 16.6552, 0.6894),
    @(7.8165, 0.6069, 15.0252, 0.7984, 17.2755, 0.8108, 15.8641, 0.8140),
    @(8.2655, 0.5936, 15.6766, 0.8091, 16.5129, 0.7903, 18.7441, 0.7975)
)

for ($i = 0; $i -lt $results.Length; $i++) {
    $row = 29 + $i
    $rowValues = $results[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

# Update the sheet's remembered scroll/selection position.
$ws.Range("E28").Select() | Out-Null
